# Append three new arrival rows (60-62) to the "Main Data" sheet, continuing
# the "Friday, Jan 13" block that already ends at row 59.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 60; Number = 59; Date = "Friday, Jan 13"; Time = "4:05 PM";  Flight = "LO3807"; From = "Warsaw"; Short = "(WAW)"; Airline = "LOT ";                     Model = "E190"; AircraftId = "(SP-LME)"; Status = "3:59 PM"; Difference = "0 hours, -6 minutes" },
    @{ Row = 61; Number = 60; Date = "Friday, Jan 13"; Time = "8:00 PM";  Flight = "FR2136"; From = "London"; Short = "(STN)"; Airline = "Ryanair ";                 Model = "B738"; AircraftId = "(EI-DYB)"; Status = "7:48 PM"; Difference = "0 hours, -12 minutes" },
    @{ Row = 62; Number = 61; Date = "Friday, Jan 13"; Time = "9:15 PM";  Flight = "FR5141"; From = "Dublin"; Short = "(DUB)"; Airline = "Ryanair (Boeing Livery) "; Model = "B738"; AircraftId = "(EI-DCL)"; Status = "9:05 PM"; Difference = "0 hours, -10 minutes" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Number
    $ws.Cells.Item($row, 2).Value = $r.Date
    $ws.Cells.Item($row, 3).Value = $r.Time
    $ws.Cells.Item($row, 4).Value = $r.Flight
    $ws.Cells.Item($row, 5).Value = $r.From
    $ws.Cells.Item($row, 6).Value = $r.Short
    $ws.Cells.Item($row, 7).Value = $r.Airline
    $ws.Cells.Item($row, 8).Value = $r.Model
    $ws.Cells.Item($row, 9).Value = $r.AircraftId
    $ws.Cells.Item($row, 10).Value = $r.Status
    $ws.Cells.Item($row, 12).Value = $r.Difference
}
